$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the ComposeItemId (column J) values for each weapon row
# -- this sets the item compose cost source per row
$composeItemIdByRow = @{
    4 = 22100401
    6 = 22100205
    7 = 22100406
    8 = 22100206
    9 = 22100208
    10 = 22100214
    11 = 22100218
    12 = 22100210
    13 = 22100201
    14 = 22100207
    15 = 22100213
    16 = 22100205
    17 = 22100401
    18 = 22100219
    19 = 22100406
    20 = 22100405
    21 = 22100406
    22 = 22100405
    23 = 22100416
    24 = 22100416
    25 = 22100406
    26 = 22100404
    27 = 22100410
    28 = 22100412
    31 = 22100415
    32 = 22100407
    33 = 22100407
    34 = 22100408
    35 = 22100411
    36 = 22100416
    39 = 22100202
    40 = 22100414
    41 = 22100209
    42 = 22100201
    43 = 22100202
    44 = 22100213
    45 = 22100414
    46 = 22100215
    47 = 22100205
    48 = 22100411
    49 = 22100204
    50 = 22100204
    51 = 22100215
    52 = 22100215
    53 = 22100206
    54 = 22100203
    55 = 22100203
    56 = 22100214
    57 = 22100219
    58 = 22100214
    59 = 22100218
}

foreach ($row in $composeItemIdByRow.Keys) {
    $ws.Cells.Item($row, 10).Value = $composeItemIdByRow[$row]
}

# Restore the last active selection recorded for this sheet
$ws.Range("G10").Select()

